$d = $word.ActiveDocument

$replacements = @(
    @("2023-11-12 Sunday", "2023-11-13 Monday"),
    @("29×81=", "51×73="),
    @("90×70=", "66×65="),
    @("69×48=", "93×81="),
    @("83×87=", "30×76="),
    @("23×96=", "87×13="),
    @("32×94=", "39×78="),
    @("79×79=", "20×79="),
    @("92×20=", "59×75="),
    @("96×67=", "91×83="),
    @("50×85=", "17×24="),
    @("46×63=", "99×49="),
    @("90×18=", "12×83="),
    @("67×23=", "21×32="),
    @("35×54=", "23×72="),
    @("95×20=", "92×29="),
    @("11×34=", "25×67="),
    @("95×45=", "68×80="),
    @("72×23=", "25×11="),
    @("76×77=", "69×14="),
    @("33×43=", "11×23="),
    @("62×47=", "45×15="),
    @("89×81=", "84×39="),
    @("69×52=", "94×88="),
    @("46×96=", "76×18="),
    @("42×11=", "66×85=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
